# Update the prediction results sheet with the finalized (hard) classification
# values: correctness labels in column B and one-hot-ish probability values
# in columns C:F (angry, happy, neutral, sad) for rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "correct"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

$ws.Range("C8").Value = 0.3386495073309164
$ws.Range("D8").Value = 0.3306152359259172
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.3307352567431665

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("B10").Value = "correct"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

$ws.Range("C12").Value = 0.6653830943429145
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.3346169056570855

$ws.Range("B13").Value = "correct"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = "correct"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0

$ws.Range("B15").Value = "correct"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0

$ws.Range("B16").Value = "correct"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0

$ws.Range("B17").Value = "correct"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0

$ws.Range("B18").Value = "correct"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0

$ws.Range("B19").Value = "correct"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0

$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0

$ws.Range("B21").Value = "correct"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1

$ws.Range("C22").Value = 0.6685935153207773
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.3314064846792227
$ws.Range("F22").Value = 0

$ws.Range("B23").Value = "correct"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1

$ws.Range("C24").Value = 0.6671008803626376
$ws.Range("D24").Value = 0.3328991196373624
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

$ws.Range("C25").Value = 0.671362379512858
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0.3286376204871418
$ws.Range("F25").Value = 0

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0
